$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename the four meta-analysis worksheet tabs (the hyphenated names force
#    quoting of the sheet references inside the _FilterDatabase defined
#    names, which happens automatically as part of the rename).
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("anchor_adjust_babies_MA").Name  = "anchor_adjust_meta-analysis"
$wb.Worksheets.Item("flag_priming_MA").Name          = "flag_priming_meta-analysis"
$wb.Worksheets.Item("gender_math_IAT_MA").Name       = "math_gender_IAT_meta-analysis"
$wb.Worksheets.Item("power_performance_MA").Name     = "power_performance_meta-analysis"

# ---------------------------------------------------------------------------
# 2. Update the selections that were left behind on a few sheets, and switch
#    the active tab to the "flag_priming_meta-analysis" sheet (index 2,
#    1-based) so it ends up as the workbook's active/selected sheet.
# ---------------------------------------------------------------------------
$wsAnchorAdjust = $wb.Worksheets.Item("anchor_adjust_meta-analysis")
$wsAnchorAdjust.Select()
$wsAnchorAdjust.Range("C2").Select()

$wsGenderMath = $wb.Worksheets.Item("math_gender_IAT_meta-analysis")
$wsGenderMath.Select()
$wsGenderMath.Range("F24").Select()

$wsFlagPriming = $wb.Worksheets.Item("flag_priming_meta-analysis")
$wsFlagPriming.Select()
$wsFlagPriming.Range("E23").Select()
